$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the exception description text in B14 (merged B14:B15)
$ws.Range("B14").Value = "Exceção 1`n[Stock Inválido] (Passo 4)"
$ws.Rows("14").RowHeight = 19.5

# Update the active selection to match the saved view state
$ws.Range("G8").Select()
